$d = $word.ActiveDocument

# --- Update "Curso (semestre ideal)" line ---
$null = $d.Content.Find.Execute("Curso (semestre ideal): EQD (10), EQN (12)", $true, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EQD (9), EQN (11)", 2)

# --- Phase 1: replace each original Requisitos entry (unique text) with a unique placeholder token ---
# This avoids ambiguity since several target strings already exist elsewhere in the original list.
$null = $d.Content.Find.Execute("LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ01@@", 2)
$null = $d.Content.Find.Execute("LOB1004 -  Cálculo II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ02@@", 2)
$null = $d.Content.Find.Execute("LOB1006 -  Cálculo IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ03@@", 2)
$null = $d.Content.Find.Execute("LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ04@@", 2)
$null = $d.Content.Find.Execute("LOB1011 -  Eletricidade Aplicada  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ05@@", 2)
$null = $d.Content.Find.Execute("LOB1012 -  Estatística  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ06@@", 2)
$null = $d.Content.Find.Execute("LOB1018 -  Física I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ07@@", 2)
$null = $d.Content.Find.Execute("LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ08@@", 2)
$null = $d.Content.Find.Execute("LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ09@@", 2)
$null = $d.Content.Find.Execute("LOB1039 -  Física Experimental III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ10@@", 2)
$null = $d.Content.Find.Execute("LOB1040 -  Laboratório de Eletricidade  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ11@@", 2)
$null = $d.Content.Find.Execute("LOB1052 -  Cálculo III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ12@@", 2)
$null = $d.Content.Find.Execute("LOB1053 -  Física III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ13@@", 2)
$null = $d.Content.Find.Execute("LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ14@@", 2)
$null = $d.Content.Find.Execute("LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ15@@", 2)
$null = $d.Content.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ16@@", 2)
$null = $d.Content.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ17@@", 2)
$null = $d.Content.Find.Execute("LOB1024 -  Mecânica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ18@@", 2)
$null = $d.Content.Find.Execute("LOB1036 -  Geometria Analítica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ19@@", 2)
$null = $d.Content.Find.Execute("LOB1038 -  Física Experimental I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ20@@", 2)
$null = $d.Content.Find.Execute("LOB1037 -  Àlgebra Linear  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ21@@", 2)
$null = $d.Content.Find.Execute("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ22@@", 2)
$null = $d.Content.Find.Execute("LOB1003 -  Cálculo I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ23@@", 2)
$null = $d.Content.Find.Execute("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@REQ24@@", 2)

# --- Phase 2: replace each placeholder with the final text for that position ---
$null = $d.Content.Find.Execute("@@REQ01@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ02@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4095 -  Química Geral Experimental  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ03@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ04@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ05@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4102 -  Nivelamento em Engenharia  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ06@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4103 -  Escrita Acadêmico Científica  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ07@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1006 -  Cálculo IV  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ08@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1053 -  Física III  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ09@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ10@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1003 -  Cálculo I  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ11@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1011 -  Eletricidade Aplicada  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ12@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1012 -  Estatística  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ13@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1024 -  Mecânica  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ14@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1036 -  Geometria Analítica  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ15@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1037 -  Álgebra Linear  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ16@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1038 -  Física Experimental I  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ17@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1039 -  Física Experimental III  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ18@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1052 -  Cálculo III  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ19@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ20@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1004 -  Cálculo II  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ21@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ22@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1018 -  Física I  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ23@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1019 -  Física II  (Requisito)", 2)
$null = $d.Content.Find.Execute("@@REQ24@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", 2)
